$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 132
$ws_ALC.Range("H132").Value = 21363820
$ws_ALC.Range("I132").Value = 25101466
$ws_ALC.Range("J132").Value = 5842.7144
$ws_ALC.Range("K132").Value = 75304398
$ws_ALC.Range("L132").Value = 17528.1432
$ws_ALC.Range("M132").Value = -75301868
$ws_ALC.Range("N132").Value = -22588.1432

# ALC row 137
$ws_ALC.Range("H137").Value = 4295.3
$ws_ALC.Range("I137").Value = 3201.7715
$ws_ALC.Range("K137").Value = 9605.3145
$ws_ALC.Range("M137").Value = -7055.3145

# ALC row 138
$ws_ALC.Range("H138").Value = 3147.32
$ws_ALC.Range("I138").Value = 1926.0588
$ws_ALC.Range("J138").Value = 3397.4578
$ws_ALC.Range("K138").Value = 5778.1764
$ws_ALC.Range("L138").Value = 10192.3734
$ws_ALC.Range("M138").Value = -638.1764000000003
$ws_ALC.Range("N138").Value = -20472.3734

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 5
$ws_ARM.Range("H5").Value = 278.2
$ws_ARM.Range("I5").Value = 150
$ws_ARM.Range("J5").Value = 363.66666
$ws_ARM.Range("K5").Value = 150
$ws_ARM.Range("L5").Value = 363.66666
$ws_ARM.Range("N5").Value = -587.66666
$ws_ARM.Range("M5").Value = -38

# ARM row 32
$ws_ARM.Range("H32").Value = 3230.25
$ws_ARM.Range("I32").Value = 2798.6382
$ws_ARM.Range("J32").Value = 9992.166999999999
$ws_ARM.Range("K32").Value = 2798.6382
$ws_ARM.Range("L32").Value = 9992.166999999999
$ws_ARM.Range("M32").Value = -2511.6382
$ws_ARM.Range("N32").Value = -10566.167

# ARM row 74
$ws_ARM.Range("H74").Value = 1464.1904
$ws_ARM.Range("I74").Value = 983.17645
$ws_ARM.Range("J74").Value = 3508.5
$ws_ARM.Range("K74").Value = 983.17645
$ws_ARM.Range("L74").Value = 3508.5
$ws_ARM.Range("M74").Value = -109.17645
$ws_ARM.Range("N74").Value = -5256.5

# ARM row 77
$ws_ARM.Range("H77").Value = 1464.1904
$ws_ARM.Range("I77").Value = 983.17645
$ws_ARM.Range("J77").Value = 3508.5
$ws_ARM.Range("K77").Value = 4915.882250000001
$ws_ARM.Range("L77").Value = 17542.5
$ws_ARM.Range("M77").Value = -547.8822500000006
$ws_ARM.Range("N77").Value = -26278.5

# ARM row 97
$ws_ARM.Range("H97").Value = 1805.3334
$ws_ARM.Range("I97").Value = 1366.4
$ws_ARM.Range("J97").Value = 4000
$ws_ARM.Range("K97").Value = 1366.4
$ws_ARM.Range("L97").Value = 4000
$ws_ARM.Range("M97").Value = -870.4000000000001
$ws_ARM.Range("N97").Value = -4992

# ARM row 102
$ws_ARM.Range("H102").Value = 2009.1666
$ws_ARM.Range("I102").Value = 2009
$ws_ARM.Range("K102").Value = 2009
$ws_ARM.Range("M102").Value = -387

# ARM row 122
$ws_ARM.Range("H122").Value = 3765.7273
$ws_ARM.Range("I122").Value = 1882
$ws_ARM.Range("J122").Value = 5649.4546
$ws_ARM.Range("K122").Value = 5646
$ws_ARM.Range("L122").Value = 16948.3638
$ws_ARM.Range("M122").Value = -3196
$ws_ARM.Range("N122").Value = -21848.3638

# ARM row 132
$ws_ARM.Range("H132").Value = 1769.6897
$ws_ARM.Range("I132").Value = 896.38635
$ws_ARM.Range("J132").Value = 4514.357
$ws_ARM.Range("K132").Value = 2689.15905
$ws_ARM.Range("L132").Value = 13543.071
$ws_ARM.Range("M132").Value = -159.1590500000002
$ws_ARM.Range("N132").Value = -18603.071

# ARM row 139
$ws_ARM.Range("H139").Value = 43470.59
$ws_ARM.Range("J139").Value = 43470.59
$ws_ARM.Range("L139").Value = 43470.59
$ws_ARM.Range("N139").Value = -53750.59

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 4
$ws_BSM.Range("H4").Value = 278.2
$ws_BSM.Range("I4").Value = 150
$ws_BSM.Range("J4").Value = 363.66666
$ws_BSM.Range("K4").Value = 150
$ws_BSM.Range("L4").Value = 363.66666
$ws_BSM.Range("N4").Value = -593.66666
$ws_BSM.Range("M4").Value = -35

# BSM row 8
$ws_BSM.Range("H8").Value = 2791.5557
$ws_BSM.Range("J8").Value = 4980
$ws_BSM.Range("L8").Value = 4980
$ws_BSM.Range("N8").Value = -5260

# BSM row 99
$ws_BSM.Range("H99").Value = 4104.0586
$ws_BSM.Range("I99").Value = 1200
$ws_BSM.Range("J99").Value = 4726.357
$ws_BSM.Range("K99").Value = 1200
$ws_BSM.Range("L99").Value = 4726.357
$ws_BSM.Range("M99").Value = 298
$ws_BSM.Range("N99").Value = -7722.357

# BSM row 134
$ws_BSM.Range("H134").Value = 2234.6765
$ws_BSM.Range("I134").Value = 1268.3636
$ws_BSM.Range("J134").Value = 6322.923
$ws_BSM.Range("K134").Value = 3805.0908
$ws_BSM.Range("L134").Value = 18968.769
$ws_BSM.Range("M134").Value = -1270.0908
$ws_BSM.Range("N134").Value = -24038.769

# BSM row 137
$ws_BSM.Range("H137").Value = 40328.75
$ws_BSM.Range("J137").Value = 40328.75
$ws_BSM.Range("L137").Value = 40328.75
$ws_BSM.Range("N137").Value = -50528.75

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws_CRP.Range("H31").Value = 2805.761
$ws_CRP.Range("I31").Value = 1083.2812
$ws_CRP.Range("J31").Value = 6742.857
$ws_CRP.Range("K31").Value = 1083.2812
$ws_CRP.Range("L31").Value = 6742.857
$ws_CRP.Range("M31").Value = -788.2811999999999
$ws_CRP.Range("N31").Value = -7332.857

# CRP row 34
$ws_CRP.Range("H34").Value = 2805.761
$ws_CRP.Range("I34").Value = 1083.2812
$ws_CRP.Range("J34").Value = 6742.857
$ws_CRP.Range("K34").Value = 1083.2812
$ws_CRP.Range("L34").Value = 6742.857
$ws_CRP.Range("M34").Value = -881.2811999999999
$ws_CRP.Range("N34").Value = -7146.857

# CRP row 123
$ws_CRP.Range("H123").Value = 39280
$ws_CRP.Range("J123").Value = 39280
$ws_CRP.Range("L123").Value = 39280
$ws_CRP.Range("N123").Value = -49080

# CRP row 132
$ws_CRP.Range("H132").Value = 1958.2931
$ws_CRP.Range("I132").Value = 1700.7954
$ws_CRP.Range("J132").Value = 2767.5715
$ws_CRP.Range("K132").Value = 5102.3862
$ws_CRP.Range("L132").Value = 8302.7145
$ws_CRP.Range("M132").Value = -2572.3862
$ws_CRP.Range("N132").Value = -13362.7145

# CRP row 140
$ws_CRP.Range("H140").Value = 163190
$ws_CRP.Range("J140").Value = 163190
$ws_CRP.Range("L140").Value = 163190
$ws_CRP.Range("N140").Value = -173550

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 107
$ws_CUL.Range("H107").Value = 78486.69500000001
$ws_CUL.Range("I107").Value = 440.83334
$ws_CUL.Range("J107").Value = 145383.14
$ws_CUL.Range("K107").Value = 1322.50002
$ws_CUL.Range("L107").Value = 436149.42
$ws_CUL.Range("M107").Value = 597.4999800000001
$ws_CUL.Range("N107").Value = -439989.42

# CUL row 113
$ws_CUL.Range("H113").Value = 746.32556
$ws_CUL.Range("I113").Value = 652.23334
$ws_CUL.Range("K113").Value = 1956.70002
$ws_CUL.Range("M113").Value = 213.29998

# CUL row 122
$ws_CUL.Range("H122").Value = 2503.3872
$ws_CUL.Range("I122").Value = 556.75
$ws_CUL.Range("J122").Value = 3430.3572
$ws_CUL.Range("K122").Value = 5010.75
$ws_CUL.Range("L122").Value = 30873.2148
$ws_CUL.Range("M122").Value = -2560.75
$ws_CUL.Range("N122").Value = -35773.2148

# CUL row 130
$ws_CUL.Range("H130").Value = 2465
$ws_CUL.Range("I130").Value = 1940
$ws_CUL.Range("J130").Value = 2990
$ws_CUL.Range("K130").Value = 5820
$ws_CUL.Range("L130").Value = 8970
$ws_CUL.Range("M130").Value = -800
$ws_CUL.Range("N130").Value = -19010

# CUL row 131
$ws_CUL.Range("H131").Value = 11910645
$ws_CUL.Range("I131").Value = 62526750
$ws_CUL.Range("J131").Value = 973.6177
$ws_CUL.Range("K131").Value = 187580250
$ws_CUL.Range("L131").Value = 2920.8531
$ws_CUL.Range("M131").Value = -187575210
$ws_CUL.Range("N131").Value = -13000.8531

# CUL row 136
$ws_CUL.Range("H136").Value = 3171.5386
$ws_CUL.Range("J136").Value = 3720
$ws_CUL.Range("L136").Value = 11160
$ws_CUL.Range("N136").Value = -21360

# CUL row 137
$ws_CUL.Range("H137").Value = 7070.2085
$ws_CUL.Range("I137").Value = 2754.6155
$ws_CUL.Range("J137").Value = 12170.454
$ws_CUL.Range("K137").Value = 8263.8465
$ws_CUL.Range("L137").Value = 36511.362
$ws_CUL.Range("M137").Value = -3163.8465
$ws_CUL.Range("N137").Value = -46711.362

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 52
$ws_GSM.Range("H52").Value = 35500
$ws_GSM.Range("J52").Value = 35500
$ws_GSM.Range("L52").Value = 35500
$ws_GSM.Range("N52").Value = -36018

# GSM row 58
$ws_GSM.Range("H58").Value = 28666.666
$ws_GSM.Range("J58").Value = 28666.666
$ws_GSM.Range("L58").Value = 28666.666
$ws_GSM.Range("N58").Value = -29220.666

# GSM row 80
$ws_GSM.Range("H80").Value = 20835816
$ws_GSM.Range("I80").Value = 41668664
$ws_GSM.Range("J80").Value = 2966.6667
$ws_GSM.Range("K80").Value = 41668664
$ws_GSM.Range("L80").Value = 2966.6667
$ws_GSM.Range("M80").Value = -41667666
$ws_GSM.Range("N80").Value = -4962.6667

# GSM row 83
$ws_GSM.Range("H83").Value = 20835816
$ws_GSM.Range("I83").Value = 41668664
$ws_GSM.Range("J83").Value = 2966.6667
$ws_GSM.Range("K83").Value = 208343320
$ws_GSM.Range("L83").Value = 14833.3335
$ws_GSM.Range("M83").Value = -208338328
$ws_GSM.Range("N83").Value = -24817.3335

# GSM row 132
$ws_GSM.Range("H132").Value = 1831.9193
$ws_GSM.Range("I132").Value = 613.4838999999999
$ws_GSM.Range("J132").Value = 3050.3547
$ws_GSM.Range("K132").Value = 1840.4517
$ws_GSM.Range("L132").Value = 9151.0641
$ws_GSM.Range("M132").Value = 689.5483000000002
$ws_GSM.Range("N132").Value = -14211.0641

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 6665.1665
$ws_LTW.Range("I7").Value = 4700
$ws_LTW.Range("J7").Value = 7647.75
$ws_LTW.Range("K7").Value = 4700
$ws_LTW.Range("L7").Value = 7647.75
$ws_LTW.Range("M7").Value = -4588
$ws_LTW.Range("N7").Value = -7871.75

# LTW row 16
$ws_LTW.Range("H16").Value = 518.2593000000001
$ws_LTW.Range("I16").Value = 527.8461
$ws_LTW.Range("K16").Value = 527.8461
$ws_LTW.Range("M16").Value = -357.8461

# LTW row 40
$ws_LTW.Range("H40").Value = 5026.255
$ws_LTW.Range("I40").Value = 4718.512
$ws_LTW.Range("J40").Value = 6288
$ws_LTW.Range("K40").Value = 4718.512
$ws_LTW.Range("L40").Value = 6288
$ws_LTW.Range("M40").Value = -4582.512
$ws_LTW.Range("N40").Value = -6560

# LTW row 100
$ws_LTW.Range("H100").Value = 2072.4285
$ws_LTW.Range("I100").Value = 1917.8334
$ws_LTW.Range("J100").Value = 3000
$ws_LTW.Range("K100").Value = 1917.8334
$ws_LTW.Range("L100").Value = 3000
$ws_LTW.Range("M100").Value = -1376.8334
$ws_LTW.Range("N100").Value = -4082

# LTW row 126
$ws_LTW.Range("H126").Value = 6665.1665
$ws_LTW.Range("I126").Value = 4700
$ws_LTW.Range("J126").Value = 7647.75
$ws_LTW.Range("K126").Value = 14100
$ws_LTW.Range("L126").Value = 22943.25
$ws_LTW.Range("M126").Value = -11630
$ws_LTW.Range("N126").Value = -27883.25

# LTW row 132
$ws_LTW.Range("H132").Value = 6030.1143
$ws_LTW.Range("I132").Value = 2174.6667
$ws_LTW.Range("K132").Value = 6524.000100000001
$ws_LTW.Range("M132").Value = -3994.000100000001

# LTW row 139
$ws_LTW.Range("H139").Value = 45355
$ws_LTW.Range("J139").Value = 50710
$ws_LTW.Range("L139").Value = 50710
$ws_LTW.Range("N139").Value = -60990

# LTW row 141
$ws_LTW.Range("H141").Value = 42075.453
$ws_LTW.Range("J141").Value = 42075.453
$ws_LTW.Range("L141").Value = 42075.453
$ws_LTW.Range("N141").Value = -52435.453

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 4
$ws_WVR.Range("H4").Value = 4000
$ws_WVR.Range("I4").Value = 0
$ws_WVR.Range("J4").Value = 4000
$ws_WVR.Range("K4").Value = 0
$ws_WVR.Range("L4").Value = 4000
$ws_WVR.Range("M4").ClearContents()
$ws_WVR.Range("N4").Value = -4226

# WVR row 5
$ws_WVR.Range("H5").Value = 100805360
$ws_WVR.Range("I5").Value = 334338340
$ws_WVR.Range("J5").Value = 719800.1
$ws_WVR.Range("K5").Value = 334338340
$ws_WVR.Range("L5").Value = 719800.1
$ws_WVR.Range("M5").Value = -334338228
$ws_WVR.Range("N5").Value = -720024.1

# WVR row 107
$ws_WVR.Range("H107").Value = 1091.3636
$ws_WVR.Range("I107").Value = 1091.3636
$ws_WVR.Range("J107").Value = 0
$ws_WVR.Range("K107").Value = 3274.0908
$ws_WVR.Range("L107").Value = 0
$ws_WVR.Range("M107").Value = -1354.0908
$ws_WVR.Range("N107").ClearContents()

# WVR row 136
$ws_WVR.Range("H136").Value = 4115.8335
$ws_WVR.Range("I136").Value = 1149.0667
$ws_WVR.Range("J136").Value = 9060.444
$ws_WVR.Range("K136").Value = 3447.2001
$ws_WVR.Range("L136").Value = 27181.332
$ws_WVR.Range("M136").Value = -897.2001
$ws_WVR.Range("N136").Value = -32281.332
